$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25
$ws.Cells.Item(25, 1).Value = "Ice Strength"
$ws.Cells.Item(25, 2).Value = "Overtime you will gain 5% in strength for a total of +15% to your strength at level 3 and 8% in durability for a total of +24% durability at level 3. You will also gain 5% in attack and 10% in armour class (defence) for a total of +15% attack and +30% in armour class at level 3."
$ws.Cells.Item(25, 3).Value = 0.05
$ws.Cells.Item(25, 5).Value = 0.08
$ws.Cells.Item(25, 10).Value = 0.05
$ws.Cells.Item(25, 11).Value = 0.1
$ws.Cells.Item(25, 13).Value = 3
$ws.Cells.Item(25, 14).Value = 250

# Row 26
$ws.Cells.Item(26, 1).Value = "Knights Gaurdian Rose"
$ws.Cells.Item(26, 2).Value = "Overtime as you level this, you will gain 8% strength and 10% durability for a total of +24% strengh and +30% durability at level 3. You will also gain 15% armour class (defence) over time for a total of +45% armour class (defence) at level 3."
$ws.Cells.Item(26, 3).Value = 0.08
$ws.Cells.Item(26, 5).Value = 0.1
$ws.Cells.Item(26, 11).Value = 0.15
$ws.Cells.Item(26, 13).Value = 3
$ws.Cells.Item(26, 14).Value = 500
$ws.Cells.Item(26, 15).Value = "Ice Strength"
$ws.Cells.Item(26, 16).Value = 3

# Row 27
$ws.Cells.Item(27, 1).Value = "Barbarians Frozen Rage"
$ws.Cells.Item(27, 2).Value = "Over time gain +8% strength and +10% attack for a total of +32% strength and +40% attack when level 4."
$ws.Cells.Item(27, 3).Value = 0.08
$ws.Cells.Item(27, 10).Value = 0.1
$ws.Cells.Item(27, 13).Value = 4
$ws.Cells.Item(27, 14).Value = 1000
$ws.Cells.Item(27, 15).Value = "Ice Strength"
$ws.Cells.Item(27, 16).Value = 3

# Row 28
$ws.Cells.Item(28, 1).Value = "Shadows Strength"
$ws.Cells.Item(28, 2).Value = "Gain additional strength over time at 4% for a total of +20% at level 5"
$ws.Cells.Item(28, 3).Value = 0.04
$ws.Cells.Item(28, 13).Value = 5
$ws.Cells.Item(28, 14).Value = 600
$ws.Cells.Item(28, 15).Value = "Ice Strength"
$ws.Cells.Item(28, 16).Value = 2

# Row 29
$ws.Cells.Item(29, 1).Value = "Whispering Death"
$ws.Cells.Item(29, 2).Value = "Gain +3% strength and +5% durability and 8% armour class (defence) for a total of +18% strength, +30% durability and +48% armour class (defence) at level 6"
$ws.Cells.Item(29, 3).Value = 0.3
$ws.Cells.Item(29, 5).Value = 0.5
$ws.Cells.Item(29, 11).Value = 0.08
$ws.Cells.Item(29, 13).Value = 6
$ws.Cells.Item(29, 14).Value = 1000
$ws.Cells.Item(29, 15).Value = "Shadows Strength"
$ws.Cells.Item(29, 16).Value = 4

# Row 30
$ws.Cells.Item(30, 1).Value = "Knights Honor"
$ws.Cells.Item(30, 2).Value = "Overtime gain +10% armour class (defence) and +10% attack for a total of +40% attack and armour class (defence) at level 4"
$ws.Cells.Item(30, 10).Value = 0.1
$ws.Cells.Item(30, 11).Value = 0.1
$ws.Cells.Item(30, 13).Value = 4
$ws.Cells.Item(30, 14).Value = 800
$ws.Cells.Item(30, 15).Value = "Knights Gaurdian Rose"
$ws.Cells.Item(30, 16).Value = 3

# Row 31
$ws.Cells.Item(31, 1).Value = "Ice Armour"
$ws.Cells.Item(31, 2).Value = "Gain +12% armour class (defence) and +12% durabaility over time for a total of +36% armour class (defence) and +36% durability at level 3"
$ws.Cells.Item(31, 5).Value = 0.12
$ws.Cells.Item(31, 11).Value = 0.12
$ws.Cells.Item(31, 13).Value = 3
$ws.Cells.Item(31, 14).Value = 1000
$ws.Cells.Item(31, 15).Value = "Barbarians Frozen Rage"
$ws.Cells.Item(31, 16).Value = 3

# Row 32
$ws.Cells.Item(32, 1).Value = "Crown of Roses"
$ws.Cells.Item(32, 2).Value = "Gain +5% strength, durability and +10% attack over time for a toal of +30% strength, durability and +60% attack at level 6"
$ws.Cells.Item(32, 3).Value = 0.05
$ws.Cells.Item(32, 5).Value = 0.05
$ws.Cells.Item(32, 10).Value = 0.1
$ws.Cells.Item(32, 13).Value = 6
$ws.Cells.Item(32, 14).Value = 1200
$ws.Cells.Item(32, 15).Value = "Knights Honor"
$ws.Cells.Item(32, 16).Value = 3

# Row 33
$ws.Cells.Item(33, 1).Value = "Icey Slash"
$ws.Cells.Item(33, 2).Value = "Gain +6% attack over time for a total of +30% attack at level 5"
$ws.Cells.Item(33, 10).Value = 0.06
$ws.Cells.Item(33, 13).Value = 5
$ws.Cells.Item(33, 14).Value = 1000
$ws.Cells.Item(33, 15).Value = "Ice Armour"
$ws.Cells.Item(33, 16).Value = 3

# Column B width (best-fit autofit side effect of longer description text)
$ws.Columns.Item(2).ColumnWidth = 329.3333333333333
